# Update data files - Bot run at 2026-02-19 10:51:47 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 corresponds to @code2careerai
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "2026-02-19T10:51:31.451853+00:00"
$ws.Range("H8").Value = 9
$ws.Range("L8").Value = "[67735, 67733, 67737, 67734, 67742, 67743, 67746, 67753, 67749]"
